# Update "想去人数" (want-to-go count) figures on the 展览 and 全部类型 sheets
# to reflect newly generated output (gh-pages data refresh at 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F2").Value = 36
$wsExpo.Range("F3").Value = 775
$wsExpo.Range("F8").Value = 3776
$wsExpo.Range("F9").Value = 77
$wsExpo.Range("F10").Value = 4448
$wsExpo.Range("F12").Value = 1111

# --- Sheet "全部类型" ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 36
$wsAll.Range("F3").Value = 775
$wsAll.Range("F9").Value = 3776
$wsAll.Range("F10").Value = 77
$wsAll.Range("F11").Value = 4448
$wsAll.Range("F13").Value = 1111
